$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing wording on the "About" and "Logout" dropdown requirement rows
$ws.Range("B16").Value = "About element should be clickable, and take you to https://saucelabs.com/."
$ws.Range("B17").Value = "The Logout element should be clickable, and log you out from the home page and return you to the login page."

# Insert a new row for the "Reset App State" requirement (3.1.4), right after 3.1.3,
# copying the formatting (fill/border/font/row height) of the row above it.
$ws.Rows.Item(18).Insert()
$ws.Range("A17:N17").Copy()
$ws.Range("A18:N18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = 35.1
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "3.1.4"
$ws.Range("B18").Value = "The ResetAppState element should be clickable, and reset the app to its default state."

# Restore the sheet view so nothing looks scrolled/selected oddly
$ws.Range("D28").Select()
